$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the "Price" column whose new values look numeric need the
# cell format set to Text first, otherwise Excel auto-converts the
# assigned string into a number (and would drop e.g. trailing zeros).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "26.292.93"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "1.621.21"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "212.15"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "0.487"
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("D9").Value = "0.0615"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").Value = "18.85"
$ws.Range("E10").Value = "  +5.06%  "
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").Value = "1.847.46"
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("D13").Value = "1.623.51"
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("D14").Value = "4.01"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "26.306.43"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "62.51"
$ws.Range("E17").Value = "  +4.15%  "
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D20").Value = "202.00"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "4.29"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").Value = "9.35"
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("D23").Value = "6.05"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D24").Value = "1.91"
$ws.Range("E24").Value = "  +5.75%  "
$ws.Range("D25").Value = "144.01"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").Value = "15.19"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").Value = "6.56"
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("D30").Value = "0.0526"
$ws.Range("E30").Value = "  +10.83%  "
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("D32").Value = "3.18"
$ws.Range("E32").Value = "  +2.00%  "
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("D36").Value = "1.179.73"
$ws.Range("E36").Value = "  +4.82%  "
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").Value = "0.808"
$ws.Range("E38").Value = "  +2.98%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").Value = "0.496"
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("D42").Value = "0.788"
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("D43").Value = "5.34"
$ws.Range("E43").Value = "  +5.12%  "
$ws.Range("D44").Value = "1.759.14"
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("D45").Value = "93.45"
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("E46").Value = "  +14.35%  "
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("D48").Value = "53.89"
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("E51").Value = "  -0.27%  "

Write-Output "done"
